$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (A1:D1) to new English column names
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# Title-case "de/del/el/y" -> "De/Del/El/Y" in state/municipality names, and fix one value
$ws.Range('B8').Value = 'Pabellón De Arteaga'
$ws.Range('B9').Value = 'Rincón De Romos'
$ws.Range('B10').Value = 'San Francisco De Los Romo'
$ws.Range('B33').Value = 'Chiapa De Corzo'
$ws.Range('D57').Value = 0.009741902834008095
$ws.Range('B82').Value = 'Guadalupe Y Calvo'
$ws.Range('B85').Value = 'Hidalgo Del Parral'
$ws.Range('B107').Value = 'San Francisco De Borja'
$ws.Range('B108').Value = 'San Francisco De Conchos'
$ws.Range('B109').Value = 'San Francisco Del Oro'
$ws.Range('B115').Value = 'Valle De Zaragoza'
$ws.Range('B130').Value = 'San Juan De Sabinas'
$ws.Range('A141').Value = 'Ciudad De México'
$ws.Range('B145').Value = 'Cuajimalpa De Morelos'
$ws.Range('B159').Value = 'Coneto De Comonfort'
$ws.Range('B173').Value = 'Nombre De Dios'
$ws.Range('B176').Value = 'Pánuco De Coronado'
$ws.Range('B183').Value = 'San Juan De Guadalupe'
$ws.Range('B184').Value = 'San Juan Del Río'
$ws.Range('B185').Value = 'San Luis Del Cordero'
$ws.Range('B186').Value = 'San Pedro Del Gallo'
$ws.Range('A194').Value = 'Estado De México'
$ws.Range('B194').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B198').Value = 'Atizapán De Zaragoza'
$ws.Range('B205').Value = 'Coacalco De Berriozábal'
$ws.Range('B210').Value = 'Ecatepec De Morelos'
$ws.Range('B218').Value = 'Naucalpan De Juárez'
$ws.Range('B222').Value = 'San Felipe Del Progreso'
$ws.Range('B223').Value = 'San Martín De Las Pirámides'
$ws.Range('B233').Value = 'Tlalnepantla De Baz'
$ws.Range('B237').Value = 'Valle De Bravo'
$ws.Range('B238').Value = 'Villa De Allende'
$ws.Range('B248').Value = 'San Miguel De Allende'
$ws.Range('B249').Value = 'Apaseo El Alto'
$ws.Range('B250').Value = 'Apaseo El Grande'
$ws.Range('B256').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B259').Value = 'Jaral Del Progreso'
$ws.Range('B269').Value = 'San Diego De La Unión'
$ws.Range('B271').Value = 'San Francisco Del Rincón'
$ws.Range('B272').Value = 'San Luis De La Paz'
$ws.Range('B273').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B274').Value = 'Silao De La Victoria'
$ws.Range('B278').Value = 'Valle De Santiago'
$ws.Range('B283').Value = 'Acapulco De Juárez'
$ws.Range('B286').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B287').Value = 'Alcozauca De Guerrero'
$ws.Range('B290').Value = 'Atenango Del Río'
$ws.Range('B291').Value = 'Atlamajalcingo Del Monte'
$ws.Range('B292').Value = 'Atoyac De Álvarez'
$ws.Range('B293').Value = 'Ayutla De Los Libres'
$ws.Range('B295').Value = 'Chilapa De Álvarez'
$ws.Range('B296').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B297').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B301').Value = 'Coyuca De Benítez'
$ws.Range('B302').Value = 'Coyuca De Catalán'
$ws.Range('B306').Value = 'Cutzamala De Pinzón'
$ws.Range('B310').Value = 'Iguala De La Independencia'
$ws.Range('B325').Value = 'Técpan De Galeana'
$ws.Range('B327').Value = 'Tixtla De Guerrero'
$ws.Range('B331').Value = 'Tlapa De Comonfort'
$ws.Range('B343').Value = 'Atotonilco El Grande'
$ws.Range('B348').Value = 'Cuautepec De Hinojosa'
$ws.Range('B356').Value = 'Jacala De Ledezma'
$ws.Range('B359').Value = 'Mineral De La Reforma'
$ws.Range('B360').Value = 'Mineral Del Chico'
$ws.Range('B361').Value = 'Mixquiahuala De Juárez'
$ws.Range('B362').Value = 'Molango De Escamilla'
$ws.Range('B364').Value = 'Pachuca De Soto'
$ws.Range('B366').Value = 'Progreso De Obregón'
$ws.Range('B371').Value = 'Tenango De Doria'
$ws.Range('B373').Value = 'Tepehuacán De Guerrero'
$ws.Range('B374').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B375').Value = 'Tezontepec De Aldama'
$ws.Range('B381').Value = 'Tula De Allende'
$ws.Range('B382').Value = 'Tulancingo De Bravo'
$ws.Range('B385').Value = 'Zacualtipán De Ángeles'
$ws.Range('B389').Value = 'Ahualulco De Mercado'
$ws.Range('B394').Value = 'Atotonilco El Alto'
$ws.Range('B402').Value = 'Encarnación De Díaz'
$ws.Range('B408').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B409').Value = 'Ixtlahuacán Del Río'
$ws.Range('B417').Value = 'Lagos De Moreno'
$ws.Range('B421').Value = 'Ojuelos De Jalisco'
$ws.Range('B424').Value = 'San Cristóbal De La Barranca'
$ws.Range('B425').Value = 'San Juan De Los Lagos'
$ws.Range('B426').Value = 'San Juanito De Escobedo'
$ws.Range('B429').Value = 'San Miguel El Alto'
$ws.Range('B431').Value = 'Talpa De Allende'
$ws.Range('B432').Value = 'Tamazula De Gordiano'
$ws.Range('B437').Value = 'Teocuitatlán De Corona'
$ws.Range('B438').Value = 'Tepatitlán De Morelos'
$ws.Range('B439').Value = 'Tizapán El Alto'
$ws.Range('B440').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B446').Value = 'Unión De San Antonio'
$ws.Range('B447').Value = 'Unión De Tula'
$ws.Range('B450').Value = 'Yahualica De González Gallo'
$ws.Range('B453').Value = 'Zapotlán Del Rey'
$ws.Range('B454').Value = 'Zapotlán El Grande'
$ws.Range('B469').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B525').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B549').Value = 'Puente De Ixtla'
$ws.Range('B551').Value = 'Tetela Del Volcán'
$ws.Range('B561').Value = 'Santa María Del Oro'
$ws.Range('B580').Value = 'San Nicolás De Los Garza'
$ws.Range('B584').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B586').Value = 'Ayoquezco De Aldama'
$ws.Range('B588').Value = 'Coicoyán De Las Flores'
$ws.Range('B590').Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Range('B592').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B593').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B594').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B596').Value = 'Magdalena Yodocono De Porfirio Díaz'
$ws.Range('B598').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B599').Value = 'Oaxaca De Juárez'
$ws.Range('B600').Value = 'Pinotepa De Don Luis'
$ws.Range('B601').Value = 'Putla Villa De Guerrero'
$ws.Range('B651').Value = 'Tanetze De Zaragoza'
$ws.Range('B652').Value = 'Teotitlán De Flores Magón'
$ws.Range('B654').Value = 'Villa De Chilapa De Díaz'
$ws.Range('B655').Value = 'Villa De Etla'
$ws.Range('B656').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B658').Value = 'Zimatlán De Álvarez'
$ws.Range('B668').Value = 'Cuetzalan Del Progreso'
$ws.Range('B675').Value = 'Izúcar De Matamoros'
$ws.Range('B680').Value = 'Palmar De Bravo'
$ws.Range('B686').Value = 'San Salvador El Verde'
$ws.Range('B688').Value = 'Tecali De Herrera'
$ws.Range('B692').Value = 'Tepanco De López'
$ws.Range('B693').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B695').Value = 'Teteles De Avila Castillo'
$ws.Range('B697').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B708').Value = 'Amealco De Bonfil'
$ws.Range('B710').Value = 'Cadereyta De Montes'
$ws.Range('B714').Value = 'Jalpan De Serra'
$ws.Range('B715').Value = 'Pinal De Amoles'
$ws.Range('B717').Value = 'San Juan Del Río'
$ws.Range('B729').Value = 'Cerro De San Pedro'
$ws.Range('B731').Value = 'Ciudad Del Maíz'
$ws.Range('B739').Value = 'Mexquitic De Carmona'
$ws.Range('B744').Value = 'San Ciro De Acosta'
$ws.Range('B747').Value = 'Santa María Del Río'
$ws.Range('B749').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B756').Value = 'Villa De Arriaga'
$ws.Range('B757').Value = 'Villa De Ramos'
$ws.Range('B758').Value = 'Villa De Reyes'
$ws.Range('B830').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B831').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B834').Value = 'Tetla De La Solidaridad'
$ws.Range('B845').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B849').Value = 'Amatlán De Los Reyes'
$ws.Range('B851').Value = 'Boca Del Río'
$ws.Range('B861').Value = 'Cosamaloapan De Carpio'
$ws.Range('B862').Value = 'Cosautlán De Carvajal'
$ws.Range('B872').Value = 'Hueyapan De Ocampo'
$ws.Range('B873').Value = 'Ignacio De La Llave'
$ws.Range('B876').Value = 'Ixhuatlán De Madero'
$ws.Range('B883').Value = 'Juchique De Ferrer'
$ws.Range('B886').Value = 'Lerdo De Tejada'
$ws.Range('B889').Value = 'Martínez De La Torre'
$ws.Range('B898').Value = 'Paso De Ovejas'
$ws.Range('B901').Value = 'Poza Rica De Hidalgo'
$ws.Range('B907').Value = 'Soledad De Doblado'
$ws.Range('B929').Value = 'Vega De Alatorre'
$ws.Range('B935').Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range('B944').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B963').Value = 'Nochistlán De Mejía'
$ws.Range('B964').Value = 'Noria De Ángeles'
$ws.Range('B974').Value = 'Teúl De González Ortega'
$ws.Range('B975').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B979').Value = 'Villa De Cos'

# Remove obsolete footer/metadata rows (988-992), leaving row 987 blank gap removed too
$ws.Rows('988:992').Delete()
